# Update BunkerPrices at 2025-04-04 13:21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Swap the "Montreal*" (AB) / "Hong Kong" (AC) columns - header and
#    all existing data rows (1-17). Use Value2 so dates/numbers round
#    trip as numerics rather than formatted text.
# ------------------------------------------------------------------
$abRange = $ws.Range("AB1:AB17")
$acRange = $ws.Range("AC1:AC17")

$abValues = $abRange.Value2
$acValues = $acRange.Value2

$abRange.Value2 = $acValues
$acRange.Value2 = $abValues

# ------------------------------------------------------------------
# 2) The last existing row (17) loses its "final row" date-only
#    formatting now that a new last row is being appended below it.
# ------------------------------------------------------------------
$ws.Range("E17").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ------------------------------------------------------------------
# 3) Append the new data row (18) with the latest bunker prices.
# ------------------------------------------------------------------
$newRow = 18
$ws.Cells.Item($newRow, 1).Value2 = 544
$ws.Cells.Item($newRow, 2).Value2 = 570
$ws.Cells.Item($newRow, 3).Value2 = 872
$ws.Cells.Item($newRow, 4).Value2 = 522
$ws.Cells.Item($newRow, 5).Value2 = 45750
$ws.Cells.Item($newRow, 5).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow, 6).Value2 = 521
$ws.Cells.Item($newRow, 7).Value2 = 524
$ws.Cells.Item($newRow, 8).Value2 = 587.25
$ws.Cells.Item($newRow, 9).Value2 = 574
$ws.Cells.Item($newRow, 10).Value2 = 550
$ws.Cells.Item($newRow, 11).Value2 = 551
$ws.Cells.Item($newRow, 12).Value2 = 592
$ws.Cells.Item($newRow, 13).Value2 = 485
$ws.Cells.Item($newRow, 14).Value2 = 595
$ws.Cells.Item($newRow, 15).Value2 = 644
$ws.Cells.Item($newRow, 16).Value2 = 523
$ws.Cells.Item($newRow, 17).Value2 = 570
$ws.Cells.Item($newRow, 18).Value2 = 524
$ws.Cells.Item($newRow, 19).Value2 = 555
$ws.Cells.Item($newRow, 20).Value2 = 656
$ws.Cells.Item($newRow, 21).Value2 = 616
$ws.Cells.Item($newRow, 22).Value2 = 588
$ws.Cells.Item($newRow, 23).Value2 = 615
$ws.Cells.Item($newRow, 24).Value2 = 521
$ws.Cells.Item($newRow, 25).Value2 = 577
$ws.Cells.Item($newRow, 26).Value2 = 766
$ws.Cells.Item($newRow, 27).Value2 = 523
$ws.Cells.Item($newRow, 28).Value2 = 538
$ws.Cells.Item($newRow, 29).Value2 = 658
$ws.Cells.Item($newRow, 30).Value2 = 621
$ws.Cells.Item($newRow, 31).Value2 = 524
$ws.Cells.Item($newRow, 32).Value2 = 550
$ws.Cells.Item($newRow, 33).Value2 = 507
$ws.Cells.Item($newRow, 34).Value2 = 655
$ws.Cells.Item($newRow, 35).Value2 = 523
$ws.Cells.Item($newRow, 36).Value2 = 539
$ws.Cells.Item($newRow, 37).Value2 = 612.5
$ws.Cells.Item($newRow, 38).Value2 = 637
$ws.Cells.Item($newRow, 39).Value2 = 530
$ws.Cells.Item($newRow, 40).Value2 = 511
$ws.Cells.Item($newRow, 41).Value2 = 570
$ws.Cells.Item($newRow, 42).Value2 = 663
$ws.Cells.Item($newRow, 43).Value2 = 510
$ws.Cells.Item($newRow, 44).Value2 = 536
$ws.Cells.Item($newRow, 45).Value2 = 485
$ws.Cells.Item($newRow, 46).Value2 = 535
$ws.Cells.Item($newRow, 47).Value2 = 776
$ws.Cells.Item($newRow, 48).Value2 = 517
